$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-15, replacing old Strike# counts
$newValues = @{
    2  = 5
    3  = 3
    4  = 4
    5  = 2
    6  = 1
    7  = 1
    8  = 3
    9  = 3
    10 = 4
    11 = 4
    12 = 0
    13 = 6
    14 = 2
    15 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
